$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1256.410888671875
$ws.Range("C2").Value = 0.9425
$ws.Range("D2").Value = 0.9204999804496765
$ws.Range("E2").Value = 1.436599969863892
$ws.Range("F2").Value = 0.6972000002861023
$ws.Range("H2").Value = 0.7431

# Row 3
$ws.Range("B3").Value = 1240.614501953125
$ws.Range("C3").Value = 0.9925
$ws.Range("D3").Value = 0.9368
$ws.Range("E3").Value = 1.758000016212463
$ws.Range("F3").Value = 0.6200000047683716
$ws.Range("H3").Value = 0.888

# Row 4
$ws.Range("B4").Value = 831.2996215820312
$ws.Range("C4").Value = 0.9711
$ws.Range("D4").Value = 0.914
$ws.Range("E4").Value = 1.813699960708618
$ws.Range("F4").Value = 0.7603999972343445
$ws.Range("H4").Value = 0.6856

# Row 5
$ws.Range("B5").Value = 812.4990844726562
$ws.Range("C5").Value = 0.8571
$ws.Range("D5").Value = 0.852
$ws.Range("E5").Value = 1.208099961280823
$ws.Range("F5").Value = 0.5893999934196472
$ws.Range("H5").Value = 0.1363

# Row 6
$ws.Range("B6").Value = 1114.16357421875
$ws.Range("C6").Value = 0.8815
$ws.Range("D6").Value = 0.8749
$ws.Range("E6").Value = 1.097100019454956
$ws.Range("F6").Value = 0.7152000069618225
$ws.Range("H6").Value = 0.3396

# Row 7
$ws.Range("B7").Value = 866.1511840820312
$ws.Range("C7").Value = 0.8723
$ws.Range("D7").Value = 0.8701000213623047
$ws.Range("E7").Value = 1.058599948883057
$ws.Range("F7").Value = 0.7325999736785889
$ws.Range("H7").Value = 0.2972

# Row 8
$ws.Range("B8").Value = 948.3499145507812
$ws.Range("C8").Value = 0.8498
$ws.Range("D8").Value = 0.8459
$ws.Range("E8").Value = 1.070799946784973
$ws.Range("F8").Value = 0.7418000102043152
$ws.Range("H8").Value = 0.083

# Row 9
$ws.Range("B9").Value = 7069.4892578125
$ws.Range("C9").Value = 0.9111
$ws.Range("D9").Value = 0.8848
$ws.Range("E9").Value = 1.813699960708618
$ws.Range("F9").Value = 0.5893999934196472
$ws.Range("H9").Value = 3.1728
